$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.608.10"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "1.635.35"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.28"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  -1.49%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.97"
$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "1.867.70"

$ws.Range("D13").Value = "1.647.16"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.49"
$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("D17").Value = "27.609.24"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.10"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("E23").Value = "  +4.28%  "

$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.43"
$ws.Range("E25").Value = "  +1.98%  "

$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").Value = "1.452.23"

$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0166"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("E40").Value = "  +6.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.03"
$ws.Range("E41").Value = "  +8.40%  "

$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("E44").Value = "  +1.61%  "

$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "1.777.37"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +2.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.30"
$ws.Range("E49").Value = "  -1.96%  "

$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("E51").Value = "  -0.65%  "
